$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Data")

# New row 21: Saturday, Jan 14 flight FR9891 to Milan (BGY), Ryanair B738,
# aircraft (9H-QEC), status 8:50 PM, difference 0 hours, 25 minutes.
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Saturday, Jan 14"
$ws.Range("C21").Value = "8:25 PM"
$ws.Range("D21").Value = "FR9891"
$ws.Range("E21").Value = "Milan"
$ws.Range("F21").Value = "(BGY)"
$ws.Range("G21").Value = "Ryanair "
$ws.Range("H21").Value = "B738"
$ws.Range("I21").Value = "(9H-QEC)"
$ws.Range("J21").Value = "8:50 PM"
$ws.Range("L21").Value = "0 hours, 25 minutes"

# K21 and M21 stay blank (as in every other row) but still need to exist as
# real cells in the sheet, matching the K/M columns used throughout the
# table. Touching the (already-default) font size materializes the cell
# without altering its style.
$ws.Range("K21").Font.Size = 11
$ws.Range("M21").Font.Size = 11
